# Scheduled runner update: refresh cached Universalis market-price columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) for a
# handful of leve rows across the per-job sheets. Columns H-N only; no
# structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 592.96295
$ws.Range("I6").Value = 577.3077
$ws.Range("K6").Value = 1731.9231
$ws.Range("M6").Value = -1619.9231

$ws.Range("H111").Value = 7293.5
$ws.Range("I111").Value = 2461.75
$ws.Range("J111").Value = 26620.5
$ws.Range("K111").Value = 7385.25
$ws.Range("L111").Value = 79861.5
$ws.Range("M111").Value = -4318.25
$ws.Range("N111").Value = -85995.5

$ws.Range("H112").Value = 2240
$ws.Range("J112").Value = 2460
$ws.Range("L112").Value = 7380
$ws.Range("N112").Value = -9596

$ws.Range("H115").Value = 9512.611000000001
$ws.Range("I115").Value = 711.1667
$ws.Range("J115").Value = 13913.333
$ws.Range("K115").Value = 2133.5001
$ws.Range("L115").Value = 41739.999
$ws.Range("M115").Value = -566.5001000000002
$ws.Range("N115").Value = -44873.999

$ws.Range("H118").Value = 6952.1113
$ws.Range("I118").Value = 388.57144
$ws.Range("J118").Value = 9249.35
$ws.Range("K118").Value = 1165.71432
$ws.Range("L118").Value = 27748.05
$ws.Range("M118").Value = 491.28568
$ws.Range("N118").Value = -31062.05

$ws.Range("H127").Value = 6335.478
$ws.Range("I127").Value = 12848.25
$ws.Range("J127").Value = 2862
$ws.Range("K127").Value = 38544.75
$ws.Range("L127").Value = 8586
$ws.Range("M127").Value = -33584.75
$ws.Range("N127").Value = -18506

$ws.Range("H129").Value = 628329.4
$ws.Range("I129").Value = 1617
$ws.Range("J129").Value = 1434102.4
$ws.Range("K129").Value = 4851
$ws.Range("L129").Value = 4302307.199999999
$ws.Range("M129").Value = 149
$ws.Range("N129").Value = -4312307.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 481.05
$ws.Range("I97").Value = 475.33334
$ws.Range("J97").Value = 498.2
$ws.Range("K97").Value = 475.33334
$ws.Range("L97").Value = 498.2
$ws.Range("M97").Value = 20.66665999999998
$ws.Range("N97").Value = -1490.2

$ws.Range("H110").Value = 639.4545000000001
$ws.Range("I110").Value = 536
$ws.Range("J110").Value = 820.5
$ws.Range("K110").Value = 536
$ws.Range("L110").Value = 820.5
$ws.Range("M110").Value = 1509
$ws.Range("N110").Value = -4910.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1399.7858
$ws.Range("I99").Value = 871.2857
$ws.Range("J99").Value = 1928.2858
$ws.Range("K99").Value = 871.2857
$ws.Range("L99").Value = 1928.2858
$ws.Range("M99").Value = 626.7143
$ws.Range("N99").Value = -4924.2858

$ws.Range("H105").Value = 1804.2727
$ws.Range("I105").Value = 1805
$ws.Range("J105").Value = 1801.5714
$ws.Range("K105").Value = 1805
$ws.Range("L105").Value = 1801.5714
$ws.Range("M105").Value = -58
$ws.Range("N105").Value = -5295.5714

$ws.Range("H107").Value = 1030
$ws.Range("I107").Value = 895
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 895
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 1025
$ws.Range("N107").Value = -5140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1172.7778
$ws.Range("I16").Value = 998
$ws.Range("J16").Value = 1312.6
$ws.Range("K16").Value = 998
$ws.Range("L16").Value = 1312.6
$ws.Range("M16").Value = -711
$ws.Range("N16").Value = -1886.6

$ws.Range("H31").Value = 2744428.8
$ws.Range("I31").Value = 1489362.5
$ws.Range("J31").Value = 6258614.5
$ws.Range("K31").Value = 1489362.5
$ws.Range("L31").Value = 6258614.5
$ws.Range("M31").Value = -1489067.5
$ws.Range("N31").Value = -6259204.5

$ws.Range("H34").Value = 2744428.8
$ws.Range("I34").Value = 1489362.5
$ws.Range("J34").Value = 6258614.5
$ws.Range("K34").Value = 1489362.5
$ws.Range("L34").Value = 6258614.5
$ws.Range("M34").Value = -1489160.5
$ws.Range("N34").Value = -6259018.5

$ws.Range("H113").Value = 1172.7778
$ws.Range("I113").Value = 998
$ws.Range("J113").Value = 1312.6
$ws.Range("K113").Value = 998
$ws.Range("L113").Value = 1312.6
$ws.Range("M113").Value = 1172
$ws.Range("N113").Value = -5652.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 277374.2
$ws.Range("I7").Value = 338966.22
$ws.Range("K7").Value = 1016898.66
$ws.Range("M7").Value = -1016786.66

$ws.Range("H80").Value = 6557.2856
$ws.Range("I80").Value = 7267.3335
$ws.Range("J80").Value = 6363.636
$ws.Range("K80").Value = 21802.0005
$ws.Range("L80").Value = 19090.908
$ws.Range("M80").Value = -20866.0005
$ws.Range("N80").Value = -20962.908

$ws.Range("H83").Value = 6557.2856
$ws.Range("I83").Value = 7267.3335
$ws.Range("J83").Value = 6363.636
$ws.Range("K83").Value = 65406.0015
$ws.Range("L83").Value = 57272.724
$ws.Range("M83").Value = -60726.0015
$ws.Range("N83").Value = -66632.724

$ws.Range("H92").Value = 4878673
$ws.Range("J92").Value = 6098146
$ws.Range("L92").Value = 18294438
$ws.Range("N92").Value = -18296934

$ws.Range("H122").Value = 1541.625
$ws.Range("I122").Value = 407.875
$ws.Range("J122").Value = 2675.375
$ws.Range("K122").Value = 3670.875
$ws.Range("L122").Value = 24078.375
$ws.Range("M122").Value = -1220.875
$ws.Range("N122").Value = -28978.375

$ws.Range("H131").Value = 25344.146
$ws.Range("I131").Value = 496.45456
$ws.Range("J131").Value = 34454.965
$ws.Range("K131").Value = 1489.36368
$ws.Range("L131").Value = 103364.895
$ws.Range("M131").Value = 3550.63632
$ws.Range("N131").Value = -113444.895

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1283.6666
$ws.Range("I61").Value = 1150.5
$ws.Range("J61").Value = 1550
$ws.Range("K61").Value = 1150.5
$ws.Range("L61").Value = 1550
$ws.Range("M61").Value = -948.5
$ws.Range("N61").Value = -1954

$ws.Range("H113").Value = 1283.6666
$ws.Range("I113").Value = 1150.5
$ws.Range("J113").Value = 1550
$ws.Range("K113").Value = 1150.5
$ws.Range("L113").Value = 1550
$ws.Range("M113").Value = 1019.5
$ws.Range("N113").Value = -5652.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 301.86957
$ws.Range("I113").Value = 197.08333
$ws.Range("J113").Value = 416.18182
$ws.Range("K113").Value = 591.24999
$ws.Range("L113").Value = 1248.54546
$ws.Range("M113").Value = 1578.75001
$ws.Range("N113").Value = -5588.54546
